$wb = $excel.ActiveWorkbook

# The workbook currently has two sheets named "Clientes" (1st) and "Planilha1" (2nd).
# We need to swap their names: the 1st sheet becomes "Planilha1" and the 2nd becomes "Clientes".
# Use a temporary name to avoid a naming collision while swapping.

$wsClientes = $wb.Worksheets.Item("Clientes")
$wsPlanilha1 = $wb.Worksheets.Item("Planilha1")

$wsClientes.Name = "__TempRename__"
$wsPlanilha1.Name = "Clientes"
$wsClientes.Name = "Planilha1"
